# Enhance Siege Analytics "Data Engineering and Infrastructure Architecture"
# section with three new bullet points describing voter-file and boundary
# estimation work, inserted right after the section's sub-header paragraph.

$d = $word.ActiveDocument

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Data Engineering and Infrastructure Architecture") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find anchor paragraph 'Data Engineering and Infrastructure Architecture'"
}

$newBullets = @(
    "• Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections",
    "• Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government",
    "• Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations"
)

$cur = $anchor
foreach ($bullet in $newBullets) {
    $cur.Range.InsertParagraphAfter()
    $cur = $cur.Next()
    $cur.Range.Text = $bullet
}
